$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2067914.9
$ws.Range("J17").Value = 2067914.9
$ws.Range("L17").Value = 6203744.699999999
$ws.Range("N17").Value = -6204080.699999999

$ws.Range("H39").Value = 172.76471
$ws.Range("I39").Value = 79.13333
$ws.Range("J39").Value = 875
$ws.Range("K39").Value = 237.39999
$ws.Range("L39").Value = 2625
$ws.Range("M39").Value = 58.60001
$ws.Range("N39").Value = -3217

$ws.Range("H53").Value = 3243.7856
$ws.Range("I53").Value = 4652.778
$ws.Range("J53").Value = 707.6
$ws.Range("K53").Value = 4652.778
$ws.Range("L53").Value = 707.6
$ws.Range("M53").Value = -4015.778
$ws.Range("N53").Value = -1981.6

$ws.Range("H64").Value = 58663.547
$ws.Range("I64").Value = 102467.336
$ws.Range("J64").Value = 6099
$ws.Range("K64").Value = 102467.336
$ws.Range("L64").Value = 6099
$ws.Range("M64").Value = -102219.336
$ws.Range("N64").Value = -6595

$ws.Range("H67").Value = 58663.547
$ws.Range("I67").Value = 102467.336
$ws.Range("J67").Value = 6099
$ws.Range("K67").Value = 102467.336
$ws.Range("L67").Value = 6099
$ws.Range("M67").Value = -101609.336
$ws.Range("N67").Value = -7815

$ws.Range("H98").Value = 26199.088
$ws.Range("I98").Value = 32011.611
$ws.Range("J98").Value = 5274
$ws.Range("K98").Value = 32011.611
$ws.Range("L98").Value = 5274
$ws.Range("M98").Value = -30513.611
$ws.Range("N98").Value = -8270

$ws.Range("H100").Value = 69255.82000000001
$ws.Range("I100").Value = 72834.57000000001
$ws.Range("J100").Value = 62993
$ws.Range("K100").Value = 72834.57000000001
$ws.Range("L100").Value = 62993
$ws.Range("M100").Value = -72293.57000000001
$ws.Range("N100").Value = -64075

$ws.Range("H122").Value = 26199.088
$ws.Range("I122").Value = 32011.611
$ws.Range("J122").Value = 5274
$ws.Range("K122").Value = 96034.833
$ws.Range("L122").Value = 15822
$ws.Range("M122").Value = -93584.833
$ws.Range("N122").Value = -20722

$ws.Range("H133").Value = 96037
$ws.Range("J133").Value = 96037
$ws.Range("L133").Value = 96037
$ws.Range("N133").Value = -106157

$ws.Range("H138").Value = 3530.1711
$ws.Range("J138").Value = 4096.719
$ws.Range("L138").Value = 12290.157
$ws.Range("N138").Value = -22570.157

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 200
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H32").Value = 11322.667
$ws.Range("I32").Value = 10752.404
$ws.Range("K32").Value = 10752.404
$ws.Range("M32").Value = -10465.404

$ws.Range("H61").Value = 7900.488
$ws.Range("I61").Value = 8788.286
$ws.Range("J61").Value = 5988.3076
$ws.Range("K61").Value = 8788.286
$ws.Range("L61").Value = 5988.3076
$ws.Range("M61").Value = -8576.286
$ws.Range("N61").Value = -6412.3076

$ws.Range("H74").Value = 1665.1875
$ws.Range("I74").Value = 606
$ws.Range("J74").Value = 2599.7646
$ws.Range("K74").Value = 606
$ws.Range("L74").Value = 2599.7646
$ws.Range("M74").Value = 268
$ws.Range("N74").Value = -4347.7646

$ws.Range("H77").Value = 1665.1875
$ws.Range("I77").Value = 606
$ws.Range("J77").Value = 2599.7646
$ws.Range("K77").Value = 3030
$ws.Range("L77").Value = 12998.823
$ws.Range("M77").Value = 1338
$ws.Range("N77").Value = -21734.823

$ws.Range("H97").Value = 8338838
$ws.Range("I97").Value = 7332.25
$ws.Range("J97").Value = 25001850
$ws.Range("K97").Value = 7332.25
$ws.Range("L97").Value = 25001850
$ws.Range("M97").Value = -6836.25
$ws.Range("N97").Value = -25002842

$ws.Range("H110").Value = 2826.3125
$ws.Range("I110").Value = 2382.5454
$ws.Range("K110").Value = 2382.5454
$ws.Range("M110").Value = -337.5454

$ws.Range("H136").Value = 7900.488
$ws.Range("I136").Value = 8788.286
$ws.Range("J136").Value = 5988.3076
$ws.Range("K136").Value = 26364.858
$ws.Range("L136").Value = 17964.9228
$ws.Range("M136").Value = -23814.858
$ws.Range("N136").Value = -23064.9228

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5797.1396
$ws.Range("I31").Value = 6016.7715
$ws.Range("J31").Value = 4836.25
$ws.Range("K31").Value = 6016.7715
$ws.Range("L31").Value = 4836.25
$ws.Range("M31").Value = -5721.7715
$ws.Range("N31").Value = -5426.25

$ws.Range("H34").Value = 5797.1396
$ws.Range("I34").Value = 6016.7715
$ws.Range("J34").Value = 4836.25
$ws.Range("K34").Value = 6016.7715
$ws.Range("L34").Value = 4836.25
$ws.Range("M34").Value = -5814.7715
$ws.Range("N34").Value = -5240.25

$ws.Range("H97").Value = 64000
$ws.Range("J97").Value = 64000
$ws.Range("L97").Value = 64000
$ws.Range("N97").Value = -65982

$ws.Range("H118").Value = 51675
$ws.Range("J118").Value = 51675
$ws.Range("L118").Value = 51675
$ws.Range("N118").Value = -54989

$ws.Range("H132").Value = 1925.75
$ws.Range("I132").Value = 1720.8334
$ws.Range("K132").Value = 5162.5002
$ws.Range("M132").Value = -2632.5002

$ws.Range("H141").Value = 178554.36
$ws.Range("J141").Value = 188633.39
$ws.Range("L141").Value = 188633.39
$ws.Range("N141").Value = -198993.39

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 3561.25
$ws.Range("I48").Value = 3357.1428
$ws.Range("J48").Value = 4990
$ws.Range("K48").Value = 10071.4284
$ws.Range("L48").Value = 14970
$ws.Range("M48").Value = -9821.428400000001
$ws.Range("N48").Value = -15470

$ws.Range("H63").Value = 2328
$ws.Range("I63").Value = 1995.6666
$ws.Range("K63").Value = 5986.9998
$ws.Range("M63").Value = -5237.9998

$ws.Range("H66").Value = 2328
$ws.Range("I66").Value = 1995.6666
$ws.Range("K66").Value = 17960.9994
$ws.Range("M66").Value = -14216.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 12659.8125
$ws.Range("I122").Value = 9514.272000000001
$ws.Range("J122").Value = 19580
$ws.Range("K122").Value = 28542.816
$ws.Range("L122").Value = 58740
$ws.Range("M122").Value = -26092.816
$ws.Range("N122").Value = -63640

$ws.Range("H123").Value = 22710.3
$ws.Range("J123").Value = 22710.3
$ws.Range("L123").Value = 22710.3
$ws.Range("N123").Value = -27610.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1434.4667
$ws.Range("I82").Value = 1243.4546
$ws.Range("J82").Value = 1959.75
$ws.Range("K82").Value = 1243.4546
$ws.Range("L82").Value = 1959.75
$ws.Range("M82").Value = -882.4546
$ws.Range("N82").Value = -2681.75

$ws.Range("H85").Value = 1434.4667
$ws.Range("I85").Value = 1243.4546
$ws.Range("J85").Value = 1959.75
$ws.Range("K85").Value = 1243.4546
$ws.Range("L85").Value = 1959.75
$ws.Range("M85").Value = 4.545399999999972
$ws.Range("N85").Value = -4455.75

$ws.Range("H122").Value = 6046.1055
$ws.Range("I122").Value = 5634.7144
$ws.Range("J122").Value = 7198
$ws.Range("K122").Value = 16904.1432
$ws.Range("L122").Value = 21594
$ws.Range("M122").Value = -14454.1432
$ws.Range("N122").Value = -26494

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 387615.7
$ws.Range("I62").Value = 950001
$ws.Range("J62").Value = 12692.167
$ws.Range("K62").Value = 950001
$ws.Range("L62").Value = 12692.167
$ws.Range("M62").Value = -949377
$ws.Range("N62").Value = -13940.167

$ws.Range("H65").Value = 387615.7
$ws.Range("I65").Value = 950001
$ws.Range("J65").Value = 12692.167
$ws.Range("K65").Value = 4750005
$ws.Range("L65").Value = 63460.835
$ws.Range("M65").Value = -4746885
$ws.Range("N65").Value = -69700.83499999999

$ws.Range("H107").Value = 35799.668
$ws.Range("I107").Value = 3049.5
$ws.Range("K107").Value = 9148.5
$ws.Range("M107").Value = -7228.5

$ws.Range("H126").Value = 26538.578
$ws.Range("I126").Value = 32522.2
$ws.Range("K126").Value = 97566.60000000001
$ws.Range("M126").Value = -95096.60000000001

$ws.Range("H132").Value = 9325
$ws.Range("J132").Value = 4878.1113
$ws.Range("L132").Value = 14634.3339
$ws.Range("N132").Value = -19694.3339

$ws.Range("H136").Value = 281331.1
$ws.Range("I136").Value = 286429.84
$ws.Range("K136").Value = 859289.52
$ws.Range("M136").Value = -856739.52
